$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.721.93"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.174.88"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.12"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "622.41"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.15"
$ws.Range("E7").Value = "  +6.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.373"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.172.63"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.750"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.207"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.57"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.54"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.673.46"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.748.11"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.174.92"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.72"
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.47"
$ws.Range("E20").Value = "  +9.80%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  +6.72%  "
$ws.Range("B22").Value = "PEPE"
$ws.Range("C22").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000211"
$ws.Range("E22").Value = "  -6.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "446.22"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.16"
$ws.Range("E25").Value = "  +7.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.54"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.323.64"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +46.56%  "
$ws.Range("E31").Value = "  +7.10%  "
$ws.Range("E32").Value = "  +16.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.41"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.76"
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.63"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.886"
$ws.Range("E37").Value = "  -11.15%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "513.03"
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.35"
$ws.Range("E40").Value = "  +5.25%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.454"
$ws.Range("E41").Value = "  +11.66%  "
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.87"
$ws.Range("E42").Value = "  +12.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.47"
$ws.Range("E43").Value = "  -9.84%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.16"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.719"
$ws.Range("E46").Value = "  +3.70%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "157.21"
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  +4.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0334"
$ws.Range("E50").Value = "  +14.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.48"
$ws.Range("E51").Value = "  +1.01%  "
